$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N (14), pushing the old
# "Late" / "Outstanding" columns one to the right. Excel copies the
# formatting of the column to the left for the newly inserted column.
$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Make "Repayment schedule" the active sheet and put the selection on
# S9 (previously G9, shifted right because of the new column), which
# also updates tabSelected/activeTab bookkeeping for the workbook.
$ws.Range("S9").Select()
